# Scene.xlsx: renumber village scene names so the trailing protocol-body
# (the FilePath-less "SceneName"/"SceneShowName" pair) never collapses to
# an empty/zero-length suffix - each row's name now carries an explicit,
# incrementing index, and a brand new row-16 scene ("villageScene6") is
# introduced so the body length is always > 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = "villageScene1"
$ws.Range("C11").Value = "villageScene1"

$ws.Range("B12").Value = "villageScene2"
$ws.Range("C12").Value = "villageScene2"

$ws.Range("B13").Value = "villageScene3"
$ws.Range("C13").Value = "villageScene3"

$ws.Range("B14").Value = "villageScene4"
$ws.Range("C14").Value = "villageScene4"

$ws.Range("B15").Value = "villageScene5"
$ws.Range("C15").Value = "villageScene5"

$ws.Range("B16").Value = "villageScene6"
$ws.Range("C16").Value = "villageScene6"

# Row 9 reverts to the sheet's standard height (no more forced wrap height);
# AutoFit drops the explicit height override entirely.
$ws.Rows.Item(9).AutoFit()

# Row 10 keeps an explicit (slightly shorter) custom height.
$ws.Rows.Item(10).RowHeight = 54.4

# Scroll/selection bookkeeping: no more pinned left column, and the last
# selected cell moved to D18.
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D18").Select()
